$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting the existing rows 46-54 down to 47-55.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly price observation.
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C46").Value = "Ñuble"
$ws.Range("D46").Value = 44900
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 100112037
$ws.Range("G46").Value = "Cebollín"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 240
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = 650
$ws.Range("N46").Value = "$/paquete 6 unidades"
$ws.Range("O46").Value = "Provincia de Diguillín"
$ws.Range("P46").Value = 108
$ws.Range("Q46").Value = 6
$ws.Range("R46").Value = "Hortaliza"
